$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill the new attendance-mark columns (Q and R) for every student row (3-21)
# with the same "p" mark already used across the sheet for columns E..P.
$ws.Range("Q3:R21").Value = "p"

# The summary/teacher row (22) only gets column R filled in (column Q stays
# empty there, matching the source edit).
$ws.Range("R22").Value = "p"

# Row 22 shrinks slightly and a new (blank) row 23 is appended right below
# the attendance table.
$ws.Rows.Item(22).RowHeight = 14.2
$ws.Rows.Item(23).RowHeight = 13.8

# Leave the selection on the new bottom-right cell, where editing ended.
$ws.Range("R23").Select()
